# Update rows 53-62 on the "Artfynd" sheet.
# The underlying edit re-orders the species records that occupy rows 53-62
# (columns A, B, D, E, F, G, H, Q, R travel together as one record), while
# the shared per-visit columns (C, I, P, S..AY) stay put. Rows 59/60 also
# swap which one of them carries the "larv/nymf" age-stage note (column K)
# for the "Zilora ferruginea" record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowVals($Row, $A, $B, $D, $E, $F, $G, $H, $Q, $R) {
    $ws.Range("A$Row").Value = $A
    $ws.Range("B$Row").Value = $B
    $ws.Range("D$Row").Value = $D
    $ws.Range("E$Row").Value = $E
    $ws.Range("F$Row").Value = $F
    $ws.Range("G$Row").Value = $G
    $ws.Range("H$Row").Value = $H
    $ws.Range("Q$Row").Value = $Q
    $ws.Range("R$Row").Value = $R
}

Set-RowVals 53 112230603 78633  "LC" 6456   "Skinnlav"            "Leptogium saturninum"          "(Dicks.) Nyl."                  572018 6697738
Set-RowVals 54 112230614 78633  "LC" 6456   "Skinnlav"            "Leptogium saturninum"          "(Dicks.) Nyl."                  571792 6697651
Set-RowVals 55 112230611 4711   "LC" 100299 "Thomsons trägnagare" "Cacotemnus thomsoni"           "(Kraatz, 1881)"                 571834 6697641
Set-RowVals 56 112230606 56575  "NT" 103021 "Talltita"            "Poecile montanus"              "(Conrad von Baldenstein, 1827)" 571961 6697705
Set-RowVals 57 112230604 102166 "LC" 222412 "Tibast"              "Daphne mezereum"                "L."                             571996 6697876
Set-RowVals 58 112230613 89539  "NT" 1202   "Ullticka"            "Phellinidium ferrugineofuscum" "(P.Karst.) Fiasson & Niemelä"   571799 6697620
Set-RowVals 59 112230612 12274  "NT" 102016 "Gropig brunbagge"    "Zilora ferruginea"              "(Paykull, 1798)"                571800 6697623
Set-RowVals 60 112230610 90466  "LC" 4769   "Svavelriska"         "Lactarius scrobiculatus"       "(Scop.:Fr.) Fr."                571853 6697760
Set-RowVals 61 112230608 99850  "LC" 221235 "Vårärt"              "Lathyrus vernus"                "(L.) Bernh."                    571931 6697694
Set-RowVals 62 112230605 99850  "LC" 221235 "Vårärt"              "Lathyrus vernus"                "(L.) Bernh."                    571995 6697876

# The "larv/nymf" age-stage note (and its surrounding Enhet/Kön/Aktivitet/
# Metod/Bestämningsmetod cells, all otherwise blank) belonged to the
# "Zilora ferruginea" record, which now sits in row 59 instead of row 60.
$ws.Range("J60").ClearContents()
$ws.Range("K60").ClearContents()
$ws.Range("L60").ClearContents()
$ws.Range("M60").ClearContents()
$ws.Range("N60").ClearContents()
$ws.Range("AF60").ClearContents()
$ws.Range("K59").Value = "larv/nymf"
